$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "68.111.68"
$cell.Style = $origStyle

$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.48%  "
$cell.Style = $origStyle

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.532.68"
$cell.Style = $origStyle

$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.46%  "
$cell.Style = $origStyle

$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.02%  "
$cell.Style = $origStyle

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "602.26"
$cell.Style = $origStyle

$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.45%  "
$cell.Style = $origStyle

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "184.76"
$cell.Style = $origStyle

$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +6.29%  "
$cell.Style = $origStyle

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle

$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.01%  "
$cell.Style = $origStyle

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.598"
$cell.Style = $origStyle

$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.38%  "
$cell.Style = $origStyle

$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.05%  "
$cell.Style = $origStyle

$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.70%  "
$cell.Style = $origStyle

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.448"
$cell.Style = $origStyle

$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.53%  "
$cell.Style = $origStyle

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.146.83"
$cell.Style = $origStyle

$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.51%  "
$cell.Style = $origStyle

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "32.61"
$cell.Style = $origStyle

$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +12.32%  "
$cell.Style = $origStyle

$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.24%  "
$cell.Style = $origStyle

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "68.066.60"
$cell.Style = $origStyle

$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.43%  "
$cell.Style = $origStyle

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0000182"
$cell.Style = $origStyle

$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.79%  "
$cell.Style = $origStyle

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.535.38"
$cell.Style = $origStyle

$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.72%  "
$cell.Style = $origStyle

$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.57%  "
$cell.Style = $origStyle

$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.32%  "
$cell.Style = $origStyle

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "400.52"
$cell.Style = $origStyle

$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.50%  "
$cell.Style = $origStyle

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.16"
$cell.Style = $origStyle

$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.85%  "
$cell.Style = $origStyle

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "74.00"
$cell.Style = $origStyle

$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.28%  "
$cell.Style = $origStyle

$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.58%  "
$cell.Style = $origStyle

$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.09%  "
$cell.Style = $origStyle

$cell = $ws.Range("B25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "LEO"
$cell.Style = $origStyle

$cell = $ws.Range("C25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$cell.Style = $origStyle

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.70"
$cell.Style = $origStyle

$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.19%  "
$cell.Style = $origStyle

$cell = $ws.Range("B26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "PEPE"
$cell.Style = $origStyle

$cell = $ws.Range("C26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell.Style = $origStyle

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0000125"
$cell.Style = $origStyle

$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.69%  "
$cell.Style = $origStyle

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.70"
$cell.Style = $origStyle

$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.36%  "
$cell.Style = $origStyle

$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.85%  "
$cell.Style = $origStyle

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = $origStyle

$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.07%  "
$cell.Style = $origStyle

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.34"
$cell.Style = $origStyle

$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.97%  "
$cell.Style = $origStyle

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.48"
$cell.Style = $origStyle

$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.09%  "
$cell.Style = $origStyle

$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.23%  "
$cell.Style = $origStyle

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.47"
$cell.Style = $origStyle

$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.02%  "
$cell.Style = $origStyle

$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.95%  "
$cell.Style = $origStyle

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "163.75"
$cell.Style = $origStyle

$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.53%  "
$cell.Style = $origStyle

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.96"
$cell.Style = $origStyle

$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.28%  "
$cell.Style = $origStyle

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.883"
$cell.Style = $origStyle

$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.45%  "
$cell.Style = $origStyle

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.16"
$cell.Style = $origStyle

$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.60%  "
$cell.Style = $origStyle

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.82"
$cell.Style = $origStyle

$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +7.16%  "
$cell.Style = $origStyle

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.35"
$cell.Style = $origStyle

$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.30%  "
$cell.Style = $origStyle

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.78"
$cell.Style = $origStyle

$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.91%  "
$cell.Style = $origStyle

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.917.53"
$cell.Style = $origStyle

$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.06%  "
$cell.Style = $origStyle

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.67"
$cell.Style = $origStyle

$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.15%  "
$cell.Style = $origStyle

$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.Style = $origStyle

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "42.58"
$cell.Style = $origStyle

$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.78%  "
$cell.Style = $origStyle

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "353.84"
$cell.Style = $origStyle

$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.15%  "
$cell.Style = $origStyle

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0308"
$cell.Style = $origStyle

$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.97%  "
$cell.Style = $origStyle

$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.08%  "
$cell.Style = $origStyle

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "33.77"
$cell.Style = $origStyle

$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.35%  "
$cell.Style = $origStyle

